$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Music"
$ws.Range("B5").Value = "Seth Harmon"
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = "Rendition of Test"

$ws.Range("D5").Select()
